$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for rows 2-6 (A: group id, B: count)
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 149

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 148

$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 133

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 120

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 80

# Remove rows 7-11 which no longer exist in the target sheet
$ws.Range("A7:B11").EntireRow.Delete()
